$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.839.11"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "1.779.94"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "2.035.21"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("E13").Value = "  +5.11%  "
$ws.Range("D14").Value = "1.804.99"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "33.880.98"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "238.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.82%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("E31").Value = "  -2.76%  "
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "1.387.69"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.66%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "78.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  -3.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +14.57%  "
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").Value = "0.0₆0140"
$ws.Range("E45").Value = "  +12.64%  "
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.47%  "
$ws.Range("D50").Value = "1.936.81"
$ws.Range("E50").Value = "  -1.11%  "
